$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 "H" row updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 225
$wsOff.Range("C2").Value = 156
$wsOff.Range("D2").Value = 55
$wsOff.Range("E2").Value = 32
$wsOff.Range("F2").Value = 3

# DEF sheet - Week 17 "H" row updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 213
$wsDef.Range("C2").Value = 139
$wsDef.Range("D2").Value = 59
$wsDef.Range("E2").Value = 33
